# Update "想去人数" (F column) figures for the newly generated output.
# Values were bumped for the 展览 (Exhibition) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row -> new F value for 展览
$sheet1Updates = @{
    3  = 593
    5  = 301
    6  = 1123
    8  = 587
    9  = 119
    10 = 767
    11 = 73
    12 = 188
    14 = 468
    15 = 1414
    20 = 95
    21 = 671
    22 = 1023
    24 = 260
    26 = 6089
    31 = 14819
    32 = 1470
    33 = 246
    34 = 113
    36 = 10832
    37 = 675
    38 = 4242
    39 = 177
    40 = 366
    41 = 114
}

foreach ($row in $sheet1Updates.Keys) {
    $sheet1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new F value for 全部类型
$sheet4Updates = @{
    3  = 593
    5  = 301
    6  = 1124
    8  = 587
    9  = 119
    10 = 767
    11 = 73
    12 = 188
    14 = 468
    15 = 1414
    21 = 95
    22 = 671
    24 = 1023
    26 = 260
    29 = 6089
    34 = 14819
    35 = 1470
    36 = 246
    37 = 113
    39 = 10832
    40 = 675
    41 = 4242
    42 = 177
    43 = 366
    44 = 114
}

foreach ($row in $sheet4Updates.Keys) {
    $sheet4.Range("F$row").Value = $sheet4Updates[$row]
}
